# New weekly price record is inserted at row 82 ("Fruta / hortaliza, semanal"),
# pushing all subsequent records (old rows 82-145) down by one row
# (new rows 83-146). Excel's row-insert semantics take care of the shift
# (values, shared formatting, etc.) for us; we only need to populate the
# freshly-inserted row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("82:82").Insert()

$ws.Range("A82").Value = 4
$ws.Range("B82").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C82").Value = "Los Lagos"
$ws.Range("D82").Value = 45072
$ws.Range("D82").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = 100112031
$ws.Range("G82").Value = "Poroto verde"
$ws.Range("H82").Value = "Magnum"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 40
$ws.Range("K82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = 30000
$ws.Range("N82").Value = '$/malla 25 kilos'
$ws.Range("O82").Value = "Perú"
$ws.Range("P82").Value = 1200
$ws.Range("Q82").Value = 25
$ws.Range("R82").Value = "Hortaliza"
